# Apply updates to column F ("dSF") values on Sheet1, per the recorded diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -5
    5  = 2
    6  = -2
    7  = -2
    8  = -2
    11 = 5
    14 = -1
    18 = 0
    19 = 3
    23 = -2
    29 = -3
    32 = -2
    34 = 0
    39 = -3
    42 = -1
    50 = 2
    51 = -1
    54 = 1
    56 = 0
    64 = 1
    65 = -1
    71 = 2
    72 = -3
    76 = 6
    79 = -5
    83 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
